$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 2
$ws.Range("H2").Value = 661829.2
$ws.Range("I2").Value = 909440.1
$ws.Range("J2").Value = 1533.3334
$ws.Range("K2").Value = 909440.1
$ws.Range("L2").Value = 1533.3334
$ws.Range("M2").Value = -909327.1
$ws.Range("N2").Value = -1759.3334
# Row 15
$ws.Range("H15").Value = 1998.9459
$ws.Range("I15").Value = 1998.9459
$ws.Range("K15").Value = 5996.8377
$ws.Range("M15").Value = -5827.8377
# Row 29
$ws.Range("H29").Value = 398.33334
$ws.Range("I29").Value = 398.33334
$ws.Range("J29").Value = 0
$ws.Range("K29").Value = 1195.00002
$ws.Range("L29").Value = 0
$ws.Range("M29").Value = -914.0000199999999
$ws.Range("N29").Value = $null
# Row 38
$ws.Range("H38").Value = 257.5
$ws.Range("I38").Value = 257.5
$ws.Range("K38").Value = 772.5
$ws.Range("M38").Value = -400.5
# Row 48
$ws.Range("H48").Value = 2754.4546
$ws.Range("J48").Value = 2944.3333
$ws.Range("L48").Value = 8832.999899999999
$ws.Range("N48").Value = -9416.999899999999
# Row 56
$ws.Range("H56").Value = 2754.4546
$ws.Range("J56").Value = 2944.3333
$ws.Range("L56").Value = 8832.999899999999
$ws.Range("N56").Value = -9900.999899999999
# Row 135
$ws.Range("H135").Value = 63192.812
$ws.Range("I135").Value = 513.7273
$ws.Range("J135").Value = 201086.8
$ws.Range("K135").Value = 4623.545700000001
$ws.Range("L135").Value = 1809781.2
$ws.Range("M135").Value = -2088.545700000001
$ws.Range("N135").Value = -1814851.2
# Row 137
$ws.Range("H137").Value = 1837.8422
$ws.Range("I137").Value = 1256.8823
$ws.Range("J137").Value = 2308.1428
$ws.Range("K137").Value = 3770.6469
$ws.Range("L137").Value = 6924.428400000001
$ws.Range("M137").Value = -1220.6469
$ws.Range("N137").Value = -12024.4284
# Row 138
$ws.Range("H138").Value = 6871.3447
$ws.Range("I138").Value = 7560.8335
$ws.Range("J138").Value = 6691.478
$ws.Range("K138").Value = 22682.5005
$ws.Range("L138").Value = 20074.434
$ws.Range("M138").Value = -17542.5005
$ws.Range("N138").Value = -30354.434
# Row 141
$ws.Range("H141").Value = 2285.9473
$ws.Range("I141").Value = 2002.25
$ws.Range("J141").Value = 3799
$ws.Range("K141").Value = 6006.75
$ws.Range("L141").Value = 11397
$ws.Range("M141").Value = -826.75
$ws.Range("N141").Value = -21757

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 13538.713
$ws.Range("I32").Value = 14255.276
$ws.Range("J32").Value = 8587.909
$ws.Range("K32").Value = 14255.276
$ws.Range("L32").Value = 8587.909
$ws.Range("M32").Value = -13968.276
$ws.Range("N32").Value = -9161.909
# Row 74
$ws.Range("H74").Value = 2535.2896
$ws.Range("I74").Value = 2322.2646
$ws.Range("K74").Value = 2322.2646
$ws.Range("M74").Value = -1448.2646
# Row 77
$ws.Range("H77").Value = 2535.2896
$ws.Range("I77").Value = 2322.2646
$ws.Range("K77").Value = 11611.323
$ws.Range("M77").Value = -7243.323
# Row 110
$ws.Range("H110").Value = 985.52
$ws.Range("I110").Value = 778.13635
$ws.Range("J110").Value = 2506.3333
$ws.Range("K110").Value = 778.13635
$ws.Range("L110").Value = 2506.3333
$ws.Range("M110").Value = 1266.86365
$ws.Range("N110").Value = -6596.3333

$ws = $wb.Worksheets.Item("BSM")
# Row 20
$ws.Range("H20").Value = 1413.0476
$ws.Range("I20").Value = 886.41174
$ws.Range("J20").Value = 3651.25
$ws.Range("K20").Value = 886.41174
$ws.Range("L20").Value = 3651.25
$ws.Range("M20").Value = -639.41174
$ws.Range("N20").Value = -4145.25

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 2774.3103
$ws.Range("I31").Value = 4243.625
$ws.Range("J31").Value = 2214.5715
$ws.Range("K31").Value = 4243.625
$ws.Range("L31").Value = 2214.5715
$ws.Range("M31").Value = -3948.625
$ws.Range("N31").Value = -2804.5715
# Row 34
$ws.Range("H34").Value = 2774.3103
$ws.Range("I34").Value = 4243.625
$ws.Range("J34").Value = 2214.5715
$ws.Range("K34").Value = 4243.625
$ws.Range("L34").Value = 2214.5715
$ws.Range("M34").Value = -4041.625
$ws.Range("N34").Value = -2618.5715
# Row 58
$ws.Range("H58").Value = 171400.17
$ws.Range("I58").Value = 336007
$ws.Range("J58").Value = 6793.3335
$ws.Range("K58").Value = 336007
$ws.Range("L58").Value = 6793.3335
$ws.Range("M58").Value = -335804
$ws.Range("N58").Value = -7199.3335
# Row 132
$ws.Range("H132").Value = 3445.923
$ws.Range("I132").Value = 1930.1
$ws.Range("J132").Value = 8498.666999999999
$ws.Range("K132").Value = 5790.299999999999
$ws.Range("L132").Value = 25496.001
$ws.Range("M132").Value = -3260.299999999999
$ws.Range("N132").Value = -30556.001
# Row 134
$ws.Range("H134").Value = 45001.84
$ws.Range("I134").Value = 50549.145
$ws.Range("J134").Value = 15878.5
$ws.Range("K134").Value = 151647.435
$ws.Range("L134").Value = 47635.5
$ws.Range("M134").Value = -149112.435
$ws.Range("N134").Value = -52705.5
# Row 136
$ws.Range("H136").Value = 171400.17
$ws.Range("I136").Value = 336007
$ws.Range("J136").Value = 6793.3335
$ws.Range("K136").Value = 1008021
$ws.Range("L136").Value = 20380.0005
$ws.Range("M136").Value = -1005471
$ws.Range("N136").Value = -25480.0005
# Row 141
$ws.Range("H141").Value = 0
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 0
$ws.Range("M141").Value = $null
$ws.Range("N141").Value = $null

$ws = $wb.Worksheets.Item("CUL")
# Row 61
$ws.Range("H61").Value = 149
$ws.Range("I61").Value = 149
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 447
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -232
$ws.Range("N61").Value = $null
# Row 68
$ws.Range("H68").Value = 3447.158
$ws.Range("I68").Value = 3019.8
$ws.Range("J68").Value = 3599.7856
$ws.Range("K68").Value = 9059.400000000001
$ws.Range("L68").Value = 10799.3568
$ws.Range("M68").Value = -8248.400000000001
$ws.Range("N68").Value = -12421.3568
# Row 71
$ws.Range("H71").Value = 3447.158
$ws.Range("I71").Value = 3019.8
$ws.Range("J71").Value = 3599.7856
$ws.Range("K71").Value = 27178.2
$ws.Range("L71").Value = 32398.0704
$ws.Range("M71").Value = -23122.2
$ws.Range("N71").Value = -40510.0704
# Row 112
$ws.Range("H112").Value = 9999.5
$ws.Range("J112").Value = 10750
$ws.Range("L112").Value = 32250
$ws.Range("N112").Value = -34466
# Row 129
$ws.Range("H129").Value = 7378.6
$ws.Range("I129").Value = 9609.846
$ws.Range("J129").Value = 3234.8572
$ws.Range("K129").Value = 28829.538
$ws.Range("L129").Value = 9704.571599999999
$ws.Range("M129").Value = -23829.538
$ws.Range("N129").Value = -19704.5716

$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Range("H70").Value = 4663.3
$ws.Range("I70").Value = 4516.8125
$ws.Range("J70").Value = 5249.25
$ws.Range("K70").Value = 4516.8125
$ws.Range("L70").Value = 5249.25
$ws.Range("M70").Value = -4246.8125
$ws.Range("N70").Value = -5789.25
# Row 73
$ws.Range("H73").Value = 4663.3
$ws.Range("I73").Value = 4516.8125
$ws.Range("J73").Value = 5249.25
$ws.Range("K73").Value = 4516.8125
$ws.Range("L73").Value = 5249.25
$ws.Range("M73").Value = -3580.8125
$ws.Range("N73").Value = -7121.25
# Row 97
$ws.Range("H97").Value = 740.087
$ws.Range("I97").Value = 787.6111
$ws.Range("J97").Value = 569
$ws.Range("K97").Value = 787.6111
$ws.Range("L97").Value = 569
$ws.Range("M97").Value = -291.6111
$ws.Range("N97").Value = -1561
# Row 102
$ws.Range("H102").Value = 4909.1333
$ws.Range("I102").Value = 4593.375
$ws.Range("J102").Value = 5270
$ws.Range("K102").Value = 4593.375
$ws.Range("L102").Value = 5270
$ws.Range("M102").Value = -2971.375
$ws.Range("N102").Value = -8514
# Row 113
$ws.Range("H113").Value = 105136.85
$ws.Range("I113").Value = 94232.37
$ws.Range("J113").Value = 118464.555
$ws.Range("K113").Value = 94232.37
$ws.Range("L113").Value = 118464.555
$ws.Range("M113").Value = -92062.37
$ws.Range("N113").Value = -122804.555
# Row 122
$ws.Range("H122").Value = 2309.2273
$ws.Range("I122").Value = 1694
$ws.Range("J122").Value = 3949.8333
$ws.Range("K122").Value = 5082
$ws.Range("L122").Value = 11849.4999
$ws.Range("M122").Value = -2632
$ws.Range("N122").Value = -16749.4999

$ws = $wb.Worksheets.Item("LTW")
# Row 40
$ws.Range("H40").Value = 5488.143
$ws.Range("I40").Value = 3892.6667
$ws.Range("K40").Value = 3892.6667
$ws.Range("M40").Value = -3756.6667
# Row 122
$ws.Range("H122").Value = 4609.5557
$ws.Range("I122").Value = 5248.25
$ws.Range("J122").Value = 4427.0713
$ws.Range("K122").Value = 15744.75
$ws.Range("L122").Value = 13281.2139
$ws.Range("M122").Value = -13294.75
$ws.Range("N122").Value = -18181.2139
# Row 132
$ws.Range("H132").Value = 38326.766
$ws.Range("I132").Value = 48017.81
$ws.Range("J132").Value = 6830.875
$ws.Range("K132").Value = 144053.43
$ws.Range("L132").Value = 20492.625
$ws.Range("M132").Value = -141523.43
$ws.Range("N132").Value = -25552.625

$ws = $wb.Worksheets.Item("WVR")
# Row 62
$ws.Range("H62").Value = 172407.83
$ws.Range("I62").Value = 6640
$ws.Range("J62").Value = 503943.5
$ws.Range("K62").Value = 6640
$ws.Range("L62").Value = 503943.5
$ws.Range("M62").Value = -6016
$ws.Range("N62").Value = -505191.5
# Row 65
$ws.Range("H65").Value = 172407.83
$ws.Range("I65").Value = 6640
$ws.Range("J65").Value = 503943.5
$ws.Range("K65").Value = 33200
$ws.Range("L65").Value = 2519717.5
$ws.Range("M65").Value = -30080
$ws.Range("N65").Value = -2525957.5
# Row 126
$ws.Range("H126").Value = 88871.086
$ws.Range("I126").Value = 170808.67
$ws.Range("J126").Value = 6933.5
$ws.Range("K126").Value = 512426.01
$ws.Range("L126").Value = 20800.5
$ws.Range("M126").Value = -509956.01
$ws.Range("N126").Value = -25740.5
# Row 132
$ws.Range("H132").Value = 28961.861
$ws.Range("I132").Value = 30577.266
$ws.Range("J132").Value = 1500
$ws.Range("K132").Value = 91731.798
$ws.Range("L132").Value = 4500
$ws.Range("M132").Value = -89201.798
$ws.Range("N132").Value = -9560
